# Generate Report for Handoff
# Updates the status + timestamps on the localization-status workbook
# to reflect the files having been readied / handed off.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns (B2 = zh-cn, C2 = de-de) and Latest Handoff Date (D2)
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-09-18 07:09:50"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (E2)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-18 07:09:48"

# de-de sheet: Status (C2) and Latest Handoff Datetime (E2)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-18 07:09:50"
